$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments ---
# Excel's Range/Columns.ColumnWidth (character units) maps to the raw OOXML
# <col width> as (value + 0.8333333333333333), so subtract that offset to
# land on the exact integer width required by the target file.
$offset = 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 57 - $offset   # C: 56 -> 57
$ws.Columns.Item(4).ColumnWidth = 41 - $offset   # D: 33 -> 41
$ws.Columns.Item(6).ColumnWidth = 16 - $offset   # F: 17 -> 16
$ws.Columns.Item(8).ColumnWidth = 36 - $offset   # H: 32 -> 36

# --- Data rows (A2:H17) ---
# Columns: A=Opportunity ID, B=Opportunity Link, C=Title, D=Country,
#          E=Premium, F=Applicants, G=Duration, H=Organization
$data = @(
    @("1329786","https://aiesec.org/opportunity/global-talent/1329786","Strategic Alliances Trainee","Panamá, Provincia de Panamá, Panamá","No","0 applicants","6 - 18 Months","NOVARTIS"),
    @("1329784","https://aiesec.org/opportunity/global-talent/1329784","Direct Contingent Worker Tech Service","Panamá, Provincia de Panamá, Panamá","No","0 applicants","6 - 18 Months","Alcon Centroamérica PA"),
    @("1329768","https://aiesec.org/opportunity/global-talent/1329768","Portfolio Manager Assistant","Београд, Србија","No","0 applicants","9 - 12 Weeks","Nebulaa IT Solutions"),
    @("1329657","https://aiesec.org/opportunity/global-talent/1329657","Sales and Marketing Specialist","Ankara, Türkiye","No","2 applicants","9 - 12 Weeks","Metaform"),
    @("1329656","https://aiesec.org/opportunity/global-talent/1329656","Software Developer","Ankara, Türkiye","No","5 applicants","9 - 12 Weeks","Metaform"),
    @("1329321","https://aiesec.org/opportunity/global-talent/1329321","Interior Designer","Cairo, Cairo Governorate, Egypt","No","0 applicants","9 - 12 Weeks","Khaled Elhusseiny Designs"),
    @("1328744","https://aiesec.org/opportunity/global-talent/1328744","[Remote] Front Desk – Customer Service","No location available","No","0 applicants","Remote","Aurent LLC"),
    @("1328730","https://aiesec.org/opportunity/global-talent/1328730","Marketing","Bursa, Türkiye","No","21 applicants","6 - 18 Months","Orhan Holding Au"),
    @("1328625","https://aiesec.org/opportunity/global-talent/1328625","Travel Advisory intern","Hyderabad, Telangana, India","No","1 applicant","9 - 12 Weeks","Amaavi Luxe Travels"),
    @("1328442","https://aiesec.org/opportunity/global-talent/1328442","Brand Ambassador","台灣臺北","No","83 applicants","3 - 6 Months","Din Tai Fung Restaurant Co., Ltd."),
    @("1328310","https://aiesec.org/opportunity/global-talent/1328310","Full Stack Developer","Cairo, Cairo Governorate, Egypt","No","1 applicant","9 - 12 Weeks","Flip Inverted Arts Academy"),
    @("1327775","https://aiesec.org/opportunity/global-talent/1327775","Accelerate Romania| Programming Intern","Bucharest, Romania","No","70 applicants","9 - 12 Weeks","AQUAsoft"),
    @("1327768","https://aiesec.org/opportunity/global-talent/1327768","Accelerate Romania| Business Development Intern","Bucharest, Romania","No","38 applicants","9 - 12 Weeks","AQUAsoft"),
    @("1325464","https://aiesec.org/opportunity/global-talent/1325464","Accelerate Romania|Account Manager for Foreign Markets","Bucharest, Romania","No","51 applicants","9 - 12 Weeks","Azuvioo"),
    @("1321641","https://aiesec.org/opportunity/global-talent/1321641","Business Administration","Adana, Reşatbey, Seyhan/Adana, Türkiye","No","65 applicants","9 - 12 Weeks","Özsayın Soğutma"),
    @("1307242","https://aiesec.org/opportunity/global-talent/1307242","Automotive Mechatronics","Ghaziabad, India","No","28 applicants","9 - 12 Weeks","KIET group of institutions")
)

$r = 2
foreach ($row in $data) {
    # Column A holds purely numeric-looking IDs (e.g. "1329786"). Assigning
    # that directly would auto-coerce to a number cell, unlike the source
    # file which stores it as text. Prefixing with an apostrophe forces text
    # entry (Excel's quote-prefix behaviour); the apostrophe itself is not
    # part of the stored value. Resetting Style afterwards clears the
    # quote-prefix formatting flag so no stray style index is left behind.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 1).Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r++
}
